$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# A new trading record (row r=5, date serial 46062) needs to be inserted
# above the current row 5. Shift the existing data rows 5-10 down to 6-11
# (bottom-up so we don't clobber data before it's copied), preserving each
# cell's value/style, then set the J-column "Current Price" formula
# explicitly since Copy() does not carry the formula along here.
for ($r = 10; $r -ge 5; $r--) {
    $dest = $r + 1
    $ws.Range("A$r`:G$r").Copy($ws.Range("A$dest`:G$dest"))
    $ws.Cells.Item($dest, 9).Value = $ws.Cells.Item($r, 9).Value2
    $ws.Cells.Item($dest, 10).Formula = $ws.Cells.Item($r, 10).Formula
}

# Populate the new row 5 with the new trading record
$ws.Cells.Item(5, 1).Value = 46062
$ws.Cells.Item(5, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(5, 2).Value = "NSE"
$ws.Cells.Item(5, 3).Value = "Buy"
$ws.Cells.Item(5, 4).Value = 5
$ws.Cells.Item(5, 5).Value = 779.35
$ws.Cells.Item(5, 6).Value = 3916.25
$ws.Cells.Item(5, 7).Value = "CN#252611665409"
$ws.Cells.Item(5, 9).Value = 19.5
$ws.Cells.Item(5, 10).Formula = "=Index!`$C`$2"
